# Insert a new weekly record at row 151 ("Apio" / Primera, fecha 2022-01-27)
# for "Vega Monumental Concepción". Inserting the row shifts every existing
# record (rows 151-196) down by one position, which matches the diff: the
# whole block of historic rows cascades down by one row and a brand-new row
# (197) appears holding what used to be the last row (196) of data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(151).Insert()

$ws.Cells.Item(151, 1).Value  = 11
$ws.Cells.Item(151, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(151, 3).Value  = "Bíobío"
$ws.Cells.Item(151, 4).Value  = 44588
$ws.Cells.Item(151, 5).Value  = 8
$ws.Cells.Item(151, 6).Value  = 100112017
$ws.Cells.Item(151, 7).Value  = "Apio"
$ws.Cells.Item(151, 8).Value  = "Americana (o)"
$ws.Cells.Item(151, 9).Value  = "Primera"
$ws.Cells.Item(151, 10).Value = 250
$ws.Cells.Item(151, 11).Value = 6500
$ws.Cells.Item(151, 12).Value = 7000
$ws.Cells.Item(151, 13).Value = 6800
$ws.Cells.Item(151, 14).Value = "$/docena de matas"
$ws.Cells.Item(151, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(151, 16).Value = 1133
$ws.Cells.Item(151, 17).Value = 6
$ws.Cells.Item(151, 18).Value = "Hortaliza"
